$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Total" column (X) is introduced: header in row 1 ...
$ws.Cells.Item(1, 24).Value = "Total"

# ... plus the row-wise totals for the existing disease-category rows (2-7).
$ws.Cells.Item(2, 24).Value = 2
$ws.Cells.Item(3, 24).Value = 2088
$ws.Cells.Item(4, 24).Value = 257
$ws.Cells.Item(5, 24).Value = 710
$ws.Cells.Item(6, 24).Value = 379
$ws.Cells.Item(7, 24).Value = 1489

# New category row 8 ("Outros" / Others), with age-group counts plus the row total.
$outros = @("Outros", 113, 3, 6, 45, 75, 77, 78, 105, 126, 149, 204, 222, 287, 342, 346, 366, 358, 342, 294, 100, 30, 1, 3669)
for ($col = 0; $col -lt $outros.Length; $col++) {
    $ws.Cells.Item(8, $col + 1).Value = $outros[$col]
}

# New grand-"Total" row 9: column-wise sums across all category rows (2-8).
$total = @("Total", 123, 5, 6, 52, 88, 100, 119, 166, 211, 295, 418, 540, 748, 893, 892, 950, 1018, 907, 737, 268, 57, 1, 8594)
for ($col = 0; $col -lt $total.Length; $col++) {
    $ws.Cells.Item(9, $col + 1).Value = $total[$col]
}
